$d = $word.ActiveDocument
$wmain = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark from its original location (the "Add the
#    new .cpp to the list of .cpp in makefile" bullet). It will be re-added
#    below in its new home paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Turn the empty paragraph right after the intro sentence into the new
#    "The propagate_invariant.pl script" section, followed by a new (empty)
#    paragraph that just carries the relocated "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$introRng = $d.Content
$introRng.Find.Execute("This document contains steps associated with some procedures for modifying the")
$introPara = $introRng.Paragraphs(1)
$emptyPara = $introPara.Next()

$fragNewSection = '<w:p xmlns:w="' + $wmain + '">' + `
    '<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>The propagate_invariant.pl script</w:t></w:r>' + `
    '<w:r><w:br/></w:r>' + `
    '<w:r><w:br/><w:t xml:space="preserve">Run this to propagate updated definitive version of source files found in model directories with updated versions.  The script </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>soruce</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> code contains a list of definitive files to propagate, and the model containing the definitive version.</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p xmlns:w="' + $wmain + '"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$emptyPara.Range.InsertXML($fragNewSection) | Out-Null

# ---------------------------------------------------------------------------
# 3) Move <w:lastRenderedPageBreak/> from the start of the "Continue with
#    implementation..." run to the start of the "Do trial compilation..."
#    run.
# ---------------------------------------------------------------------------
$trialRng = $d.Content
$trialRng.Find.Execute("Do trial compilation of the implementation file (TableDimensionSymbol.cpp), using Ctrl-F7")
$trialPara = $trialRng.Paragraphs(1)

$fragTrial = '<w:p xmlns:w="' + $wmain + '">' + `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:lastRenderedPageBreak/><w:t>Do trial compilation of the implementation file (TableDimensionSymbol.cpp), using Ctrl-F7</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">, and fix errors as </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:t>required.</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'

$trialPara.Range.InsertXML($fragTrial) | Out-Null

$continueRng = $d.Content
$continueRng.Find.Execute("Continue with implementation of functionality for the new class.")
$continuePara = $continueRng.Paragraphs(1)

$fragContinue = '<w:p xmlns:w="' + $wmain + '">' + `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>Continue with implementation of functionality for the new class.</w:t></w:r>' + `
    '<w:r><w:br/></w:r>' + `
    '<w:r><w:br/></w:r>' + `
    '</w:p>'

$continuePara.Range.InsertXML($fragContinue) | Out-Null
